# Trade #33 closed at 2026-02-17 04:15:54 - unknown UNKNOWN +0.000%

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Summary sheet: refresh aggregate stats after the new trade closed
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B5").Value = 0.41    # Total P&L %
$summary.Range("B6").Value = 33      # Total Trades
$summary.Range("B9").Value = 45.45   # Win Rate %

# ---------------------------------------------------------------------------
# 2. Strategy Status sheet: MarketMaking row picks up the new trade count
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 33       # Trades
$status.Range("G4").Value = 45.45    # Win Rate %

# ---------------------------------------------------------------------------
# Helper: append the newly closed trade's row (#34) to a trade-log sheet
# ---------------------------------------------------------------------------
function Add-TradeRow34($sheet) {
    $sheet.Cells.Item(34, 1).Value = 33

    # Date/time-like text must be forced to stay literal text (otherwise
    # Excel auto-converts "2026-02-17" into a date serial number). Apply a
    # text number format, set the value, then restore the default "Normal"
    # style so the cell doesn't retain an extra style record.
    $sheet.Cells.Item(34, 2).NumberFormat = "@"
    $sheet.Cells.Item(34, 2).Value = "2026-02-17"
    $sheet.Cells.Item(34, 2).Style = "Normal"

    $sheet.Cells.Item(34, 3).Value = "04:15:48"
    $sheet.Cells.Item(34, 4).Value = "MarketMaking"
    $sheet.Cells.Item(34, 5).Value = "DOWN"
    $sheet.Cells.Item(34, 6).Value = 0.01
    $sheet.Cells.Item(34, 7).Value = 0.01
    $sheet.Cells.Item(34, 8).Value = "CLOSED"
    $sheet.Cells.Item(34, 9).Value = 0
    $sheet.Cells.Item(34, 10).Value = 0
    $sheet.Cells.Item(34, 11).Value = 100.67
    $sheet.Cells.Item(34, 12).Value = 0
    $sheet.Cells.Item(34, 13).Value = 0
    $sheet.Cells.Item(34, 14).Value = 0.6
    $sheet.Cells.Item(34, 15).Value = "Normal spread capture: 19600 bps"
    $sheet.Cells.Item(34, 16).Value = "early_exit"
    $sheet.Cells.Item(34, 17).Value = 0.11
}

# ---------------------------------------------------------------------------
# 3. All Trades sheet: append trade #33 as row 34
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow34 $allTrades

# ---------------------------------------------------------------------------
# 4. MarketMaking sheet: append the same trade #33 as row 34
# ---------------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow34 $marketMaking
